$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("Q3").Value = 2.25
$ws.Range("R3").Value = 1.62

$ws.Range("G4").Value = 3.1
$ws.Range("H4").Value = 3.15
$ws.Range("I4").Value = 2.22
$ws.Range("J4").Value = 3.55
$ws.Range("K4").Value = 2.1
$ws.Range("L4").Value = 2.77
$ws.Range("N4").Value = 8
$ws.Range("O4").Value = 1.3
$ws.Range("P4").Value = 2.95
$ws.Range("Q4").Value = 1.93
$ws.Range("R4").Value = 1.78
$ws.Range("W4").Value = 9.75
$ws.Range("X4").Value = 17
$ws.Range("Y4").Value = 10.75
$ws.Range("Z4").Value = 40
$ws.Range("AA4").Value = 27
$ws.Range("AB4").Value = 32
$ws.Range("AC4").Value = 9.25
$ws.Range("AD4").Value = 6.1
$ws.Range("AH4").Value = 7.5
$ws.Range("AI4").Value = 10.75
$ws.Range("AJ4").Value = 8.75
$ws.Range("AK4").Value = 22
$ws.Range("AL4").Value = 18.5
$ws.Range("AM4").Value = 28
$ws.Range("AN4").Value = 5.1
$ws.Range("AO4").Value = 16.5
$ws.Range("AP4").Value = 21
$ws.Range("AQ4").Value = 75
$ws.Range("AR4").Value = 100
$ws.Range("AS4").Value = 250
$ws.Range("AT4").Value = 2.6
$ws.Range("AU4").Value = 6.6
$ws.Range("AV4").Value = 55
$ws.Range("AW4").Value = 4.15
$ws.Range("AX4").Value = 11.5
$ws.Range("AY4").Value = 18.5
$ws.Range("AZ4").Value = 45
$ws.Range("BA4").Value = 75

